$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 16: politeness_score (B16) switches from a text "3" to a real number 3. ---
# All other cells in row 16 (C16..H16) are unchanged.
$ws.Range("B16").Value = 3

# --- Row 17: brand new annotation row appended below row 16. ---
$ws.Range("A17").Value = "Sunsi Wu"

# B17 must stay a *text* "3" (like the original B16 used to be), so we can't just
# assign the bare string, because Excel auto-coerces numeric-looking strings into
# numbers. Route the value through a text formula and paste it back as a value so
# it lands as a genuine string without leaving stray number-format styles behind.
$ws.Range("ZZ1").Formula = '="3"'
$ws.Range("ZZ1").Copy()
$ws.Range("B17").PasteSpecial(-4163)
$ws.Range("ZZ1").ClearContents()

$ws.Range("C17").Value = "can"
$ws.Range("D17").Value = "SMY"
$ws.Range("E17").Value = "EXP"
$ws.Range("F17").Value = "9cb2103f-10a8-4188-b35f-b6e342d90889"
$ws.Range("G17").Value = "rJwelMbR-_annotated.xlsx"
$ws.Range("H17").Value = "The authors show through several experiments that the divide and conquer (DnC) technique can solve more complex tasks than can be solved with conventional policy gradient methods (TRPO is used as the baseline)."
